$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

$ws.Range("D1:D4").NumberFormat = "@"

$ws.Range("D1").Value = "Data podpisania kontraktu"
$ws.Range("D2").Value = "18-10-2015"
$ws.Range("D3").Value = "25-08-2015"
$ws.Range("D4").Value = "06-06-2016"

$ws.Columns.Item(4).ColumnWidth = 24.5703125

$ws.Range("D3").Select()
